# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (fund-holding detail) positioned
#    right after "2021-Q4" and before "总计".
# 2. Insert a new summary row at the top of "总计" for the 2022-Q1 quarter,
#    pushing the existing rows down and renumbering the leading index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" detail sheet right before "总计"
# ---------------------------------------------------------------------
$q4_2021 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4_2021)
$newSheet.Name = "2022-Q1"

# Header row (row 1), columns B..H - same header style used elsewhere
# (bold, centered, thin border around each cell).
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers) {
    $cell = $newSheet.Cells.Item(1, $col)
    $cell.Value = $h
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
    $col = $col + 1
}

# Data rows (A = running index, styled like the header; B..H plain data).
# Columns D, E, F, G are stored as text (matching the source data export),
# H (仓位排名) is numeric.
$rows = @(
    @("000586", "景顺长城中小板创业板精选股票", "2.42", "94.15", "5.89", "0.1425", 9),
    @("010706", "景顺长城景骊成长混合型证券投资基金", "1.13", "93.50", "5.15", "0.0582", 9),
    @("002802", "广发东财大数据精选灵活配置混合", "0.41", "55.13", "2.51", "0.0103", 4),
    @("005443", "国金量化多策略灵活配置混合", "0.51", "64.10", "0.64", "0.0033", 9)
)

$r = 2
foreach ($row in $rows) {
    $a = $newSheet.Cells.Item($r, 1)
    $a.Value = $r - 2
    $a.Font.Bold = $true
    $a.HorizontalAlignment = -4108
    $a.VerticalAlignment = -4160
    $a.Borders.LineStyle = 1

    $newSheet.Cells.Item($r, 2).NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]

    $newSheet.Cells.Item($r, 4).NumberFormat = "@"
    $newSheet.Cells.Item($r, 4).Value = $row[2]

    $newSheet.Cells.Item($r, 5).NumberFormat = "@"
    $newSheet.Cells.Item($r, 5).Value = $row[3]

    $newSheet.Cells.Item($r, 6).NumberFormat = "@"
    $newSheet.Cells.Item($r, 6).Value = $row[4]

    $newSheet.Cells.Item($r, 7).NumberFormat = "@"
    $newSheet.Cells.Item($r, 7).Value = $row[5]

    $newSheet.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: add the 2022-Q1 summary row to "总计"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push existing rows 2..6 down to 3..7, then wipe the inherited formatting
# on the freshly-inserted row 2 (Insert() copies the row-above's style).
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

$idxCell = $total.Range("A2")
$idxCell.Value = 0
$idxCell.Font.Bold = $true
$idxCell.HorizontalAlignment = -4108
$idxCell.VerticalAlignment = -4160
$idxCell.Borders.LineStyle = 1

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.21

# Renumber the leading index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
